$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue {
    param($ws, $addr, $val)
    # Write the target string via a literal-text formula, then paste-special
    # as values. This avoids Excel's "looks like a number/percent" auto
    # conversion that a direct .Value assignment would trigger for
    # numeric-looking strings like "312.06" or "-3.18%", keeping the cell a
    # plain text cell (no NumberFormat / style change) exactly like the
    # original inline-string cells.
    $escaped = $val.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextValue $ws "D2" "312.06"
Set-TextValue $ws "E2" "0.40%"
Set-TextValue $ws "D3" "38.18"
Set-TextValue $ws "E3" "-3.18%"
Set-TextValue $ws "D4" "5.137"
Set-TextValue $ws "E4" "0.31%"
Set-TextValue $ws "D5" "0.08097"
Set-TextValue $ws "E5" "-0.25%"
Set-TextValue $ws "D6" "4.448"
Set-TextValue $ws "E6" "4.96%"
Set-TextValue $ws "E7" "-2.53%"
Set-TextValue $ws "D8" "8.300"
Set-TextValue $ws "E8" "1.92%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D9" "3.263"
Set-TextValue $ws "E9" "-2.02%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D10" "0.9393"
Set-TextValue $ws "E10" "1.00%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D11" "0.1321"
Set-TextValue $ws "E11" "-7.47%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D12" "0.1955"
Set-TextValue $ws "E12" "0.92%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D13" "0.09047"
Set-TextValue $ws "E13" "-0.14%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D14" "0.03491"
Set-TextValue $ws "E14" "-0.23%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D15" "0.09692"
Set-TextValue $ws "E15" "-1.29%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D16" "0.001408"
Set-TextValue $ws "E16" "0.33%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D17" "0.005886"
Set-TextValue $ws "E17" "0.45%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D18" "3.553"
Set-TextValue $ws "E18" "-6.24%"
Set-TextValue $ws "D19" "0.3466"
Set-TextValue $ws "E19" "0.37%"
Set-TextValue $ws "D20" "0.1283"
Set-TextValue $ws "E20" "-2.24%"
Set-TextValue $ws "D21" "5.019"
Set-TextValue $ws "E21" "7.19%"
Set-TextValue $ws "E22" "2.92%"
Set-TextValue $ws "D23" "0.04370"
Set-TextValue $ws "E23" "-0.12%"
Set-TextValue $ws "D24" "0.001240"
Set-TextValue $ws "E24" "0.82%"
Set-TextValue $ws "D25" "0.004725"
Set-TextValue $ws "E25" "-1.50%"
Set-TextValue $ws "D26" "0.0003851"
Set-TextValue $ws "E26" "195.89%"
Set-TextValue $ws "D39" "0.02208"
Set-TextValue $ws "E39" "3.23%"
Set-TextValue $ws "D40" "0.05237"
Set-TextValue $ws "E40" "2.49%"
Set-TextValue $ws "D41" "0.007589"
Set-TextValue $ws "E41" "1.93%"
Set-TextValue $ws "D42" "0.01031"
Set-TextValue $ws "E42" "4.57%"
Set-TextValue $ws "D43" "0.1390"
Set-TextValue $ws "E43" "2.14%"
Set-TextValue $ws "D44" "0.002106"
Set-TextValue $ws "E44" "-1.24%"
Set-TextValue $ws "D45" "0.009111"
Set-TextValue $ws "E45" "5.58%"
Set-TextValue $ws "D46" "0.00006620"
Set-TextValue $ws "E46" "3.29%"
Set-TextValue $ws "D47" "0.00000000752"
Set-TextValue $ws "E47" "0.14%"
Set-TextValue $ws "E48" "17.99%"
Set-TextValue $ws "E49" "68.89%"
Set-TextValue $ws "D50" "0.00002106"
Set-TextValue $ws "E50" "0.14%"
Set-TextValue $ws "D51" "0.0002006"
Set-TextValue $ws "E51" "0.14%"

$excel.CutCopyMode = 0
